$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.351.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.26"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6293"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07442"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2893"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.97"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07729"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.06"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.966"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6758"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001025"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.56"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.248"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.411.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.357"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.482"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1348"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07071"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.464"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.478"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.052"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.049"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.827"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.139"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6979"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.585"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.814"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.235.16"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.83%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.797"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9301"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9996"
$ws.Range("D42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.995.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.29"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("E46").Value = "  +2.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.019"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.705"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.31%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1140"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.86%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.887"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("E51").Value = "  -0.58%  "
